$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve existing formatting on the data range, then force text
# number format so that numeric-looking strings (e.g. "1.001",
# "0.3799") are stored as text instead of being coerced to numbers,
# matching the inlineStr/text representation used in the workbook.
$dataRange = $ws.Range("D2:E51")
$origStyle = $dataRange.Style
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "23.504.86"
$ws.Range("E2").Value = "  -0.45%  "
$ws.Range("D3").Value = "1.648.87"
$ws.Range("E3").Value = "  +0.18%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.31%  "
$ws.Range("E5").Value = "  +0.33%  "
$ws.Range("E6").Value = "  -1.35%  "
$ws.Range("D7").Value = "0.3799"
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "0.3571"
$ws.Range("E8").Value = "  -0.90%  "
$ws.Range("D9").Value = "50.70"
$ws.Range("E9").Value = "  -2.57%  "
$ws.Range("D10").Value = "0.08114"
$ws.Range("E10").Value = "  -0.93%  "
$ws.Range("D11").Value = "1.224"
$ws.Range("E11").Value = "  -1.52%  "
$ws.Range("D12").Value = "1.002"
$ws.Range("E12").Value = "  +0.30%  "
$ws.Range("D13").Value = "22.12"
$ws.Range("E13").Value = "  -1.45%  "
$ws.Range("D14").Value = "6.424"
$ws.Range("E14").Value = "  -1.52%  "
$ws.Range("D15").Value = "7.422"
$ws.Range("E15").Value = "  +0.83%  "
$ws.Range("D16").Value = "0.00001205"
$ws.Range("E16").Value = "  -1.89%  "
$ws.Range("D17").Value = "1.656.23"
$ws.Range("E17").Value = "  +0.68%  "
$ws.Range("D18").Value = "97.27"
$ws.Range("E18").Value = "  +0.37%  "
$ws.Range("D19").Value = "0.06985"
$ws.Range("E19").Value = "  +0.01%  "
$ws.Range("D20").Value = "6.778"
$ws.Range("E20").Value = "  +0.56%  "
$ws.Range("D21").Value = "17.48"
$ws.Range("E21").Value = "  -0.46%  "
$ws.Range("E22").Value = "  +0.25%  "
$ws.Range("D23").Value = "12.59"
$ws.Range("E23").Value = "  +0.23%  "
$ws.Range("D24").Value = "23.518.05"
$ws.Range("E24").Value = "  -0.37%  "
$ws.Range("D25").Value = "2.482"
$ws.Range("E25").Value = "  -1.66%  "
$ws.Range("D26").Value = "2.917"
$ws.Range("E26").Value = "  -6.27%  "
$ws.Range("D27").Value = "20.99"
$ws.Range("E27").Value = "  -1.29%  "
$ws.Range("D28").Value = "152.63"
$ws.Range("E28").Value = "  +0.21%  "
$ws.Range("D29").Value = "5.230"
$ws.Range("E29").Value = "  +0.82%  "
$ws.Range("D30").Value = "133.33"
$ws.Range("D31").Value = "1.838.03"
$ws.Range("E31").Value = "  +0.42%  "
$ws.Range("D32").Value = "6.946"
$ws.Range("E32").Value = "  +2.89%  "
$ws.Range("E33").Value = "  +4.82%  "
$ws.Range("D34").Value = "11.97"
$ws.Range("E34").Value = "  +2.75%  "
$ws.Range("D35").Value = "1.032"
$ws.Range("E35").Value = "  -5.17%  "
$ws.Range("D36").Value = "0.02733"
$ws.Range("E36").Value = "  -2.14%  "
$ws.Range("D37").Value = "0.08726"
$ws.Range("E37").Value = "  -0.98%  "
$ws.Range("D38").Value = "0.2456"
$ws.Range("E38").Value = "  -2.20%  "
$ws.Range("D39").Value = "5.986"
$ws.Range("E39").Value = "  -1.48%  "
$ws.Range("D40").Value = "13.37"
$ws.Range("E40").Value = "  +4.62%  "
$ws.Range("D41").Value = "0.06876"
$ws.Range("E41").Value = "  -2.02%  "
$ws.Range("D42").Value = "0.6927"
$ws.Range("E43").Value = "  -0.27%  "
$ws.Range("D44").Value = "15.74"
$ws.Range("E44").Value = "  -0.61%  "
$ws.Range("D45").Value = "0.6454"
$ws.Range("E45").Value = "  -0.74%  "
$ws.Range("D46").Value = "1.001"
$ws.Range("E46").Value = "  +0.32%  "
$ws.Range("D47").Value = "2.274"
$ws.Range("E47").Value = "  -2.56%  "
$ws.Range("E48").Value = "  -1.25%  "
$ws.Range("D49").Value = "0.07820"
$ws.Range("E49").Value = "  -1.93%  "
$ws.Range("D50").Value = "128.35"
$ws.Range("E50").Value = "  +0.41%  "
$ws.Range("E51").Value = "  -1.33%  "

# Restore the original (default) style/format now that text values are set
$dataRange.Style = $origStyle
